# Agency_Recruits_template.xlsx
#
# Commit intent:
#   1. Show user full name rather than login name for all kinds of data.
#      -> the template placeholder for the "salesperson" column must switch
#         from ${record.salesPerson} to ${record.salesPersonFullName}.
#   2. Validate region/department required input (logic change outside the
#      scope of this template's cell content - no visible cell effect here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D2" template cell ("销售" column) currently reads ${record.salesPerson};
# switch it to the full-name placeholder. Re-setting the value naturally
# rewrites the shared-strings table (drops the old string, appends the new
# one), which is exactly the reshuffle seen across the whole sheet in the diff.
$ws.Range("D2").Value = '${record.salesPersonFullName}'

# The "销售" (salesperson) column is narrowed slightly relative to the other
# (uniform 18.5703125-wide) columns.
$ws.Columns.Item(4).ColumnWidth = 15

# Active selection moves to F6.
$ws.Range("F6").Select()
